$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 431.26666
$ws.Range("I33").Value = 305.30768
$ws.Range("J33").Value = 1250
$ws.Range("K33").Value = 305.30768
$ws.Range("L33").Value = 1250
$ws.Range("M33").Value = -76.30768
$ws.Range("N33").Value = -1708

$ws.Range("H69").Value = 15145
$ws.Range("J69").Value = 16224.75
$ws.Range("L69").Value = 48674.25
$ws.Range("N69").Value = -50422.25

$ws.Range("I70").Value = 1525225.6
$ws.Range("J70").Value = 3378.4
$ws.Range("K70").Value = 4575676.800000001
$ws.Range("L70").Value = 10135.2
$ws.Range("M70").Value = -4575406.800000001
$ws.Range("N70").Value = -10675.2

$ws.Range("H72").Value = 15145
$ws.Range("J72").Value = 16224.75
$ws.Range("L72").Value = 146022.75
$ws.Range("N72").Value = -154758.75

$ws.Range("I73").Value = 1525225.6
$ws.Range("J73").Value = 3378.4
$ws.Range("K73").Value = 4575676.800000001
$ws.Range("L73").Value = 10135.2
$ws.Range("M73").Value = -4574740.800000001
$ws.Range("N73").Value = -12007.2

$ws.Range("H104").Value = 142.5
$ws.Range("I104").Value = 93.75
$ws.Range("K104").Value = 281.25
$ws.Range("M104").Value = 1465.75

$ws.Range("H123").Value = 90999.89999999999
$ws.Range("J123").Value = 90999.89999999999
$ws.Range("L123").Value = 90999.89999999999
$ws.Range("N123").Value = -100799.9

$ws.Range("H132").Value = 6525.3335
$ws.Range("I132").Value = 2664.1428
$ws.Range("J132").Value = 13282.417
$ws.Range("K132").Value = 7992.428400000001
$ws.Range("L132").Value = 39847.251
$ws.Range("M132").Value = -5462.428400000001
$ws.Range("N132").Value = -44907.251

$ws.Range("H137").Value = 812560.7
$ws.Range("I137").Value = 797.2
$ws.Range("J137").Value = 2165499.8
$ws.Range("K137").Value = 2391.6
$ws.Range("L137").Value = 6496499.399999999
$ws.Range("M137").Value = 158.3999999999996
$ws.Range("N137").Value = -6501599.399999999


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3336.2942
$ws.Range("I45").Value = 1783
$ws.Range("J45").Value = 4423.6
$ws.Range("K45").Value = 1783
$ws.Range("L45").Value = 4423.6
$ws.Range("M45").Value = -1406
$ws.Range("N45").Value = -5177.6

$ws.Range("H102").Value = 4944.1113
$ws.Range("I102").Value = 3717.6
$ws.Range("K102").Value = 3717.6
$ws.Range("M102").Value = -2095.6

$ws.Range("H122").Value = 3632.1667
$ws.Range("I122").Value = 2540.9167
$ws.Range("K122").Value = 7622.750100000001
$ws.Range("M122").Value = -5172.750100000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2000
$ws.Range("J64").Value = 2000
$ws.Range("L64").Value = 2000
$ws.Range("N64").Value = -2450

$ws.Range("H67").Value = 2000
$ws.Range("J67").Value = 2000
$ws.Range("L67").Value = 2000
$ws.Range("N67").Value = -3560

$ws.Range("H94").Value = 2137.3235
$ws.Range("I94").Value = 1858.1482
$ws.Range("K94").Value = 1858.1482
$ws.Range("M94").Value = -1407.1482

$ws.Range("H99").Value = 2909.8333
$ws.Range("J99").Value = 4699.5
$ws.Range("L99").Value = 4699.5
$ws.Range("N99").Value = -7695.5

$ws.Range("H107").Value = 4839.4614
$ws.Range("I107").Value = 5240.1
$ws.Range("K107").Value = 5240.1
$ws.Range("M107").Value = -3320.1


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2556.5789
$ws.Range("I58").Value = 2535.0625
$ws.Range("J58").Value = 2671.3333
$ws.Range("K58").Value = 2535.0625
$ws.Range("L58").Value = 2671.3333
$ws.Range("M58").Value = -2332.0625
$ws.Range("N58").Value = -3077.3333

$ws.Range("H99").Value = 41599.2
$ws.Range("I99").Value = 14499
$ws.Range("J99").Value = 150000
$ws.Range("K99").Value = 14499
$ws.Range("L99").Value = 150000
$ws.Range("M99").Value = -13001
$ws.Range("N99").Value = -152996

$ws.Range("H126").Value = 41599.2
$ws.Range("I126").Value = 14499
$ws.Range("J126").Value = 150000
$ws.Range("K126").Value = 43497
$ws.Range("L126").Value = 450000
$ws.Range("M126").Value = -41027
$ws.Range("N126").Value = -454940

$ws.Range("H132").Value = 2884.2104
$ws.Range("I132").Value = 2794.5
$ws.Range("J132").Value = 4499
$ws.Range("K132").Value = 8383.5
$ws.Range("L132").Value = 13497
$ws.Range("M132").Value = -5853.5
$ws.Range("N132").Value = -18557

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H136").Value = 2556.5789
$ws.Range("I136").Value = 2535.0625
$ws.Range("J136").Value = 2671.3333
$ws.Range("K136").Value = 7605.1875
$ws.Range("L136").Value = 8013.999899999999
$ws.Range("M136").Value = -5055.1875
$ws.Range("N136").Value = -13113.9999


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 5490.2144
$ws.Range("I44").Value = 422.4
$ws.Range("K44").Value = 1267.2
$ws.Range("M44").Value = -869.1999999999998

$ws.Range("H68").Value = 4982.6665
$ws.Range("J68").Value = 4982.6665
$ws.Range("L68").Value = 14947.9995
$ws.Range("N68").Value = -16569.9995

$ws.Range("H71").Value = 4982.6665
$ws.Range("J71").Value = 4982.6665
$ws.Range("L71").Value = 44843.9985
$ws.Range("N71").Value = -52955.9985

$ws.Range("H130").Value = 12027.728
$ws.Range("I130").Value = 7835
$ws.Range("J130").Value = 14423.571
$ws.Range("K130").Value = 23505
$ws.Range("L130").Value = 43270.713
$ws.Range("M130").Value = -18485
$ws.Range("N130").Value = -53310.713


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 69999.5
$ws.Range("J48").Value = 69999.5
$ws.Range("L48").Value = 69999.5
$ws.Range("N48").Value = -70969.5

$ws.Range("H102").Value = 3417.6667
$ws.Range("I102").Value = 3417.6667
$ws.Range("K102").Value = 3417.6667
$ws.Range("M102").Value = -1795.6667

$ws.Range("H122").Value = 3987.2354
$ws.Range("I122").Value = 3935.5334
$ws.Range("J122").Value = 4375
$ws.Range("K122").Value = 11806.6002
$ws.Range("L122").Value = 13125
$ws.Range("M122").Value = -9356.600199999999
$ws.Range("N122").Value = -18025

$ws.Range("H126").Value = 4032.6667
$ws.Range("I126").Value = 4800
$ws.Range("J126").Value = 2498
$ws.Range("K126").Value = 14400
$ws.Range("L126").Value = 7494
$ws.Range("M126").Value = -11930
$ws.Range("N126").Value = -12434

$ws.Range("H132").Value = 10050556
$ws.Range("I132").Value = 2920.3572
$ws.Range("K132").Value = 8761.071599999999
$ws.Range("M132").Value = -6231.071599999999


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7101.2
$ws.Range("I7").Value = 7101.2
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 7101.2
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -6989.2
$ws.Range("N7").ClearContents()

$ws.Range("H40").Value = 3397.2
$ws.Range("I40").Value = 3397.2
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3397.2
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3261.2
$ws.Range("N40").ClearContents()

$ws.Range("H82").Value = 7693.625
$ws.Range("I82").Value = 4962.5
$ws.Range("J82").Value = 10424.75
$ws.Range("K82").Value = 4962.5
$ws.Range("L82").Value = 10424.75
$ws.Range("M82").Value = -4601.5
$ws.Range("N82").Value = -11146.75

$ws.Range("H85").Value = 7693.625
$ws.Range("I85").Value = 4962.5
$ws.Range("J85").Value = 10424.75
$ws.Range("K85").Value = 4962.5
$ws.Range("L85").Value = 10424.75
$ws.Range("M85").Value = -3714.5
$ws.Range("N85").Value = -12920.75

$ws.Range("H100").Value = 27808960
$ws.Range("I100").Value = 2339.6
$ws.Range("K100").Value = 2339.6
$ws.Range("M100").Value = -1798.6

$ws.Range("H122").Value = 3193.2942
$ws.Range("I122").Value = 3115.2097
$ws.Range("J122").Value = 4000.1667
$ws.Range("K122").Value = 9345.6291
$ws.Range("L122").Value = 12000.5001
$ws.Range("M122").Value = -6895.6291
$ws.Range("N122").Value = -16900.5001

$ws.Range("H126").Value = 7101.2
$ws.Range("I126").Value = 7101.2
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 21303.6
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -18833.6
$ws.Range("N126").ClearContents()


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws.Range("H48").Value = 56686
$ws.Range("J48").Value = 64999.5
$ws.Range("L48").Value = 64999.5
$ws.Range("N48").Value = -66137.5

$ws.Range("H49").Value = 102500
$ws.Range("J49").Value = 102500
$ws.Range("L49").Value = 102500
$ws.Range("N49").Value = -102960

$ws.Range("H57").Value = 200000
$ws.Range("I57").Value = 200000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 200000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -199246
$ws.Range("N57").ClearContents()

$ws.Range("H81").Value = 1849.6666
$ws.Range("I81").Value = 1248.3334
$ws.Range("K81").Value = 2496.6668
$ws.Range("M81").Value = -1435.6668

$ws.Range("H84").Value = 1849.6666
$ws.Range("I84").Value = 1248.3334
$ws.Range("K84").Value = 12483.334
$ws.Range("M84").Value = -7179.333999999999

$ws.Range("H122").Value = 4877.8096
$ws.Range("I122").Value = 4889.3335
$ws.Range("J122").Value = 4862.4443
$ws.Range("K122").Value = 14668.0005
$ws.Range("L122").Value = 14587.3329
$ws.Range("M122").Value = -12218.0005
$ws.Range("N122").Value = -19487.3329

$ws.Range("H126").Value = 11638.3
$ws.Range("I126").Value = 11825.667
$ws.Range("K126").Value = 35477.001
$ws.Range("M126").Value = -33007.001

